{"js": "// Update the TIAC template labels:\n//  1. \"Date de r\u00e9ception \u00e0 la DD(ETS)PP : {{...}}\" -> \"Date de r\u00e9ception : {{...}}\"\n//  2. \"Suite donn\u00e9e par la DD : {{...}}\"            -> \"Suite donn\u00e9e : {{...}}\"\n// (Note: the space right before the colon in these labels is a non-breaking\n//  space, U+00A0, in the original template - it is preserved here.)\n\nconst body = context.document.body;\n\nconst receptionResults = body.search(\"Date de r\u00e9ception \u00e0 la DD(ETS)PP\\u00A0: \", { matchCase: true });\nreceptionResults.load(\"items\");\nawait context.sync();\n\nif (receptionResults.items.length > 0) {\n  receptionResults.items[0].insertText(\"Date de r\u00e9ception\\u00A0: \", Word.InsertLocation.replace);\n}\n\nconst followUpResults = body.search(\"Suite donn\u00e9e par la DD\\u00A0: \", { matchCase: true });\nfollowUpResults.load(\"items\");\nawait context.sync();\n\nif (followUpResults.items.length > 0) {\n  followUpResults.items[0].insertText(\"Suite donn\u00e9e\\u00A0: \", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the TIAC template labels:\n#  1. \"Date de r\u00e9ception \u00e0 la DD(ETS)PP : {{...}}\" -> \"Date de r\u00e9ception : {{...}}\"\n#  2. \"Suite donn\u00e9e par la DD : {{...}}\"            -> \"Suite donn\u00e9e : {{...}}\"\n# (Note: the space right before the colon in these labels is a non-breaking\n#  space, U+00A0, in the original template - it is preserved here.)\n\n$d = $word.ActiveDocument\n$nbsp = [char]0x00A0\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\"Date de r\u00e9ception \u00e0 la DD(ETS)PP$($nbsp): \", $false, $false, $false, $false, $false, $true, 1, $false, \"Date de r\u00e9ception$($nbsp): \", 2)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"Suite donn\u00e9e par la DD$($nbsp): \", $false, $false, $false, $false, $false, $true, 1, $false, \"Suite donn\u00e9e$($nbsp): \", 2)\n"}
